# ex9.1.2(Linear) C_Stationary generator - "alpha_non_zero" experiment update
# Commit: "expermits todos no convexos menos el 5to"
#
# Helper: writes a value into a cell while forcing it to be stored as TEXT
# (shared-string) even when the text looks like a number - this mirrors the
# original workbook where every generated numeric value is written as a
# string. Re-applying the "Normal" cell style afterwards keeps the cell on
# the workbook's default style (so no stray per-cell style survives, even
# though the "@" number format briefly used to force text entry).
function Set-TextValue {
    param($Cell, [string]$Text)

    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# NOTE: worksheet name lookup via Worksheets.Item(<name>) is case
# insensitive, and this workbook has both "Vector_bf" and "Vector_BF"
# sheets - so every sheet is addressed by its (1-based) index to avoid
# an accidental match on the wrong one:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_follower -------------------------------------------
$ws = $wb.Worksheets.Item(3)

Set-TextValue $ws.Cells.Item(2,1) "-21.128588198390794 + 2.624669341415005y"
Set-TextValue $ws.Cells.Item(2,2) "21.128588198390794"
Set-TextValue $ws.Cells.Item(2,3) "J_0_L0_v"
Set-TextValue $ws.Cells.Item(2,4) "0.13"
Set-TextValue $ws.Cells.Item(2,5) "0"
Set-TextValue $ws.Cells.Item(2,6) "6.1"

Set-TextValue $ws.Cells.Item(3,1) "0.9044455249858698 - x + 0.645410493790575y"
Set-TextValue $ws.Cells.Item(3,2) "-3.90444552498587"
Set-TextValue $ws.Cells.Item(3,3) "J_0_L0_v"
Set-TextValue $ws.Cells.Item(3,4) "0.6"
Set-TextValue $ws.Cells.Item(3,5) "0.8999999999999999"
Set-TextValue $ws.Cells.Item(3,6) "1.5"

Set-TextValue $ws.Cells.Item(4,1) "15.218046704143411 + x - 2.648204559520921y"
Set-TextValue $ws.Cells.Item(4,2) "-27.21804670414341"
Set-TextValue $ws.Cells.Item(4,3) "J_0_LP_v"
Set-TextValue $ws.Cells.Item(4,4) "0.62"
Set-TextValue $ws.Cells.Item(4,5) "0"
Set-TextValue $ws.Cells.Item(4,6) "1.3"

Set-TextValue $ws.Cells.Item(5,1) "-24.93 + 4x"
Set-TextValue $ws.Cells.Item(5,2) "12.399999999999999"
Set-TextValue $ws.Cells.Item(5,3) "J_Ne_L0_v"
Set-TextValue $ws.Cells.Item(5,4) "0.32"
Set-TextValue $ws.Cells.Item(5,5) "5.8"
Set-TextValue $ws.Cells.Item(5,6) "0"

# --- Punto_modificado -------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Cells.Item(2,1) "6.1"
Set-TextValue $ws.Cells.Item(2,2) "8.05"

# --- Vector_bf ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Cells.Item(2,1) "-0.08656648375532472"

# --- Vector_BF -----------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Cells.Item(2,1) "-21.3"
Set-TextValue $ws.Cells.Item(3,1) "2.419130555588483"

# --- Vector_Alpha (plain numeric cell, not text) ------------------------------
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2,1).Value = 2.3241022797604605
